$wb = $excel.ActiveWorkbook

# Sheet "13" is where the new names are added (column D), and it becomes the active sheet/tab
$ws = $wb.Worksheets.Item("13")
$ws.Activate()

$ws.Range("D1").Value = "Gabreal haj"
$ws.Range("D2").Value = "Amitay Lavi "
$ws.Range("D3").Value = "Michal Yonasi "
$ws.Range("D4").Value = "Yaniv Avraham"

$ws.Range("E6").Select()
